$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new effort log entry on row 30
$ws.Range("A30").Value = 41206
$ws.Range("A30").NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Range("B30").Value = 2.75
$ws.Range("D30").Value = "Manual continued"

$ws.Range("D30").Select()
